$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Helper: replace a whole paragraph's Range contents with a literal
# OOXML <w:p> fragment (Range.InsertXML replaces the range's content).
# ---------------------------------------------------------------
function Replace-ParagraphXml($para, [string]$xml) {
    $r = $para.Range
    $r.InsertXML($xml)
}

# -----------------------------------------------------------------------
# 1) Paragraph containing the first large screenshot (cx=4751152) gets a
#    <w:lastRenderedPageBreak/> inserted right before the <w:drawing>.
# -----------------------------------------------------------------------
$frag9 = @'
<w:p xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" w14:paraId="43AA9254" w14:textId="465B7E56" w:rsidR="00EA7931" w:rsidRDefault="00C2505F" w:rsidP="00EA7931"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:lastRenderedPageBreak/><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="73514687" wp14:editId="66F7F0FD"><wp:extent cx="4751152" cy="3512820"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1527542924" name="Picture 3"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 3"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId6"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="4756336" cy="3516653"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>
'@

# -----------------------------------------------------------------------
# 2) Paragraph "We conduct monthly mock tests..." gets a
#    <w:lastRenderedPageBreak/> inserted right before its text run.
# -----------------------------------------------------------------------
$frag20 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" w14:paraId="4BF75678" w14:textId="77777777" w:rsidR="007D0A40" w:rsidRPr="00234C25" w:rsidRDefault="007D0A40" w:rsidP="007D0A40"><w:pPr><w:rPr><w:rFonts w:ascii="Oxygen" w:hAnsi="Oxygen"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="00234C25"><w:rPr><w:rFonts w:ascii="Oxygen" w:hAnsi="Oxygen"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">We conduct monthly mock tests for all popular exams across all classes and disciplines. Test schedules are available at </w:t></w:r><w:hyperlink r:id="rId8" w:history="1"><w:r w:rsidRPr="00234C25"><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="Oxygen" w:hAnsi="Oxygen"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>www.anodiam.com</w:t></w:r></w:hyperlink><w:r w:rsidRPr="00234C25"><w:rPr><w:rFonts w:ascii="Oxygen" w:hAnsi="Oxygen"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@

# -----------------------------------------------------------------------
# 3) Paragraph "Call up 9073700094 To enroll for you exam." -- the three
#    runs " To " / "enroll" (with spell-check proofErr wrapping) /
#    " for you exam." collapse into a single run's text.
# -----------------------------------------------------------------------
$frag21 = @'
<w:p xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="6ABD21AB" w14:textId="77777777" w:rsidR="007D0A40" w:rsidRPr="00234C25" w:rsidRDefault="007D0A40" w:rsidP="007D0A40"><w:pPr><w:rPr><w:rFonts w:ascii="Oxygen" w:hAnsi="Oxygen"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="00234C25"><w:rPr><w:rFonts w:ascii="Oxygen" w:hAnsi="Oxygen"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Call up </w:t></w:r><w:r w:rsidRPr="00234C25"><w:rPr><w:rFonts w:ascii="Oxygen" w:hAnsi="Oxygen"/><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>9073700094</w:t></w:r><w:r w:rsidRPr="00234C25"><w:rPr><w:rFonts w:ascii="Oxygen" w:hAnsi="Oxygen"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> To enroll for you exam.</w:t></w:r></w:p>
'@

# -----------------------------------------------------------------------
# 4) Paragraph "6. After clicking Start Test button..." loses its
#    <w:lastRenderedPageBreak/> (it moved to the picture above it).
# -----------------------------------------------------------------------
$frag32 = @'
<w:p xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="266512B7" w14:textId="0EBB669D" w:rsidR="007D0A40" w:rsidRDefault="003C0B02" w:rsidP="003C0B02"><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t xml:space="preserve">6. After clicking Start Test button, we are not able to give test because after </w:t></w:r><w:r w:rsidR="004F48CC"><w:t>clicking a</w:t></w:r><w:r><w:t xml:space="preserve"> particular topic</w:t></w:r></w:p>
'@

# -----------------------------------------------------------------------
# 5) Paragraph containing the second screenshot (cx=3756660) gets a
#    <w:lastRenderedPageBreak/> inserted right before the <w:drawing>.
# -----------------------------------------------------------------------
$frag39 = @'
<w:p xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" w14:paraId="0B2F741B" w14:textId="046E8EEE" w:rsidR="00AE6EB1" w:rsidRDefault="00AE6EB1" w:rsidP="000A3203"><w:r><w:rPr><w:noProof/></w:rPr><w:lastRenderedPageBreak/><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="2B7AEB28" wp14:editId="61E0B437"><wp:extent cx="3756660" cy="2004060"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1574157729" name="Picture 2"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 6"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId10"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="3756660" cy="2004060"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>
'@

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "*We conduct monthly mock tests*") {
        Replace-ParagraphXml $p $frag20
    } elseif ($t -like "*Call up 9073700094*") {
        Replace-ParagraphXml $p $frag21
    } elseif ($t -like "*6. After clicking Start Test button*") {
        Replace-ParagraphXml $p $frag32
    }
}

for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $s = $d.InlineShapes.Item($i)
    if ($s.Width -gt 374.0 -and $s.Width -lt 374.2 -and $s.Height -gt 276.5 -and $s.Height -lt 276.7) {
        $p = $s.Range.Paragraphs(1)
        Replace-ParagraphXml $p $frag9
    } elseif ($s.Width -gt 295.7 -and $s.Width -lt 295.9 -and $s.Height -gt 157.7 -and $s.Height -lt 157.9) {
        $p = $s.Range.Paragraphs(1)
        Replace-ParagraphXml $p $frag39
    }
}

# -----------------------------------------------------------------------
# 6) New empty ListParagraph (numId 1) inserted before the very first
#    paragraph in the document.
# -----------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$firstRange = $firstPara.Range.Duplicate
$firstRange.Collapse(1)
$firstRange.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr></w:p>")

Write-Host "done"
